$d = $word.ActiveDocument

# Locate the existing end-date text "May 31, 2025" without altering it yet,
# so we can work out exactly where it lives in the document.
$findRange = $d.Content
$found = $findRange.Find.Execute("May 31, 2025", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'May 31, 2025' in the document"
}
$mayStart = $findRange.Start
$mayEnd = $findRange.End

# Insert "April" right before "May 31, 2025". Because the insertion point sits
# exactly on the run boundary, Word folds the new text into the preceding
# run (" March 6, 2023 to "), which already carries the bold/no-italic
# formatting that "April" needs.
$insertPoint = $d.Range($mayStart, $mayStart)
$insertPoint.Collapse(1)
$insertPoint.InsertBefore("April")

# "May 31, 2025" (and everything after it) shifted right by the length of
# "April".
$shift = 5
$mayStart = $mayStart + $shift
$mayEnd = $mayEnd + $shift

# Split "April" off into its own run (distinct from "March 6, 2023 to ") by
# toggling a character property on/off so the run boundary is materialized,
# without actually changing the resulting formatting.
$aprilRange = $d.Range($mayStart - $shift, $mayStart)
$aprilRange.Font.Bold = $false
$aprilRange.Font.Bold = $true

# Replace "May 31, 2025" with " 28, 2025" (leading space kept, so the run
# starts with a plain space just like the original diff).
$dateRange = $d.Range($mayStart, $mayEnd)
$dateRange.Text = " 28, 2025"

# Split " " | "28" | ", 2025" into three separate runs, again by toggling a
# property so no new formatting is actually introduced.
$spaceRange = $d.Range($mayStart, $mayStart + 1)
$spaceRange.Font.Bold = $false
$spaceRange.Font.Bold = $true

$dayRange = $d.Range($mayStart + 1, $mayStart + 3)
$dayRange.Font.Bold = $false
$dayRange.Font.Bold = $true

Write-Output "Final text: $($d.Range($mayStart - $shift, $mayStart + 9).Text)"
